# Change the J-column "946543" numeric values into the text value "946.543"
# (PDF invoice creation feature needs this as a textual reference code, not
# a numeric amount), and update the shared contact-email column (U) to the
# new finance/ops mailbox.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jCells = @("J2", "J4", "J5")
foreach ($addr in $jCells) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = "946.543"
    $rng.ClearFormats()
}

$uRange = $ws.Range("U2:U5")
$uRange.Value = "fso@tbd-tp.bg"
